{"js": "// [Function] View Product Detail: Code View\n// Insert two new list items (\"ListProduct: Add 1\" and \"ProductDetail: Add Many\")\n// right after the \"Add product to cart\" bullet in the \"Ti\u1ebfn \u0111\u1ed9:\" (progress)\n// section, at the same list (numId 1) but one level deeper (ilvl 3).\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\n// There are two \"Add product to cart\" bullets in the document (one under the\n// \"Use cases\" outline near the top, one under \"Ti\u1ebfn \u0111\u1ed9:\" further down). The\n// new items belong right after the LAST one (the \"Ti\u1ebfn \u0111\u1ed9:\" section), so\n// keep scanning and remember the last match.\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === \"Add product to cart\") {\n    target = paragraphs.items[i];\n  }\n}\n\nif (!target) {\n  throw new Error('Could not find paragraph with text \"Add product to cart\"');\n}\n\n// Insert the first new paragraph right after the target; it inherits the\n// ListParagraph style + numbering (numId 1) from \"target\", then bump its\n// list level from 2 down to 3 to match the desired nesting.\nconst p1 = target.insertParagraph(\"ListProduct: Add 1\", Word.InsertLocation.after);\nawait context.sync();\np1.load(\"listItemOrNullObject\");\nawait context.sync();\np1.listItemOrNullObject.level = 3;\nawait context.sync();\n\n// Insert the second new paragraph right after the first one, same treatment.\nconst p2 = p1.insertParagraph(\"ProductDetail: Add Many\", Word.InsertLocation.after);\nawait context.sync();\np2.load(\"listItemOrNullObject\");\nawait context.sync();\np2.listItemOrNullObject.level = 3;\nawait context.sync();\n", "ps1": "# [Function] View Product Detail: Code View\n#\n# Insert two new bullet items (\"ListProduct: Add 1\" and \"ProductDetail: Add\n# Many\") right after the \"Add product to cart\" bullet in the \"Tien do:\"\n# (progress) section, as children one level deeper (same numbered list,\n# numId 1, but ilvl 3 instead of 2).\n\n$d = $word.ActiveDocument\n\n# The phrase \"Add product to cart\" appears twice in the document (once in\n# the \"User classes & Uses cases\" outline near the top, once further down in\n# the \"Tien do:\" / progress section). The new bullets belong after the LAST\n# occurrence, so walk all matches with Find and remember the last one.\n$rng = $d.Content\n$target = $null\nwhile ($rng.Find.Execute(\"Add product to cart\")) {\n    $target = $rng.Duplicate\n    $rng.Collapse(0)   # wdCollapseEnd - keep searching after this match\n}\n\nif ($target -eq $null) {\n    Write-Output \"Could not find paragraph 'Add product to cart'\"\n} else {\n    $targetPara = $target.Paragraphs(1)\n\n    # Insert paragraph #1 right after the target. It inherits the\n    # ListParagraph style + list numbering (numId 1, level 2) from the\n    # paragraph it follows, so bump it one level deeper (ListLevelNumber is\n    # 1-based in the COM model, so 4 == ilvl 3 in the saved XML).\n    $targetPara.Range.InsertParagraphAfter()\n    $p1 = $targetPara.Next()\n    $p1.Range.Text = \"ListProduct: Add 1\"\n    $p1.Range.ListFormat.ListLevelNumber = 4\n\n    # Insert paragraph #2 right after paragraph #1, same treatment.\n    $p1.Range.InsertParagraphAfter()\n    $p2 = $p1.Next()\n    $p2.Range.Text = \"ProductDetail: Add Many\"\n    $p2.Range.ListFormat.ListLevelNumber = 4\n\n    Write-Output \"Inserted ListProduct/ProductDetail bullets after 'Add product to cart'\"\n}\n"}
